$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.734.19"
$ws.Range("E2").Value = "  -0.27%  "

$ws.Range("D3").Value = "1.850.31"
$ws.Range("E3").Value = "  -0.96%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  -2.57%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.84"
$ws.Range("E5").Value = "  -1.47%  "

$ws.Range("E6").Value = "  -2.33%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4328"
$ws.Range("E7").Value = "  -2.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3774"
$ws.Range("E8").Value = "  -0.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07387"
$ws.Range("E9").Value = "  -1.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8850"
$ws.Range("E10").Value = "  -0.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.67"
$ws.Range("E11").Value = "  -0.64%  "

$ws.Range("D12").Value = "1.861.14"
$ws.Range("E12").Value = "  -0.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.756"
$ws.Range("E13").Value = "  -0.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.476"
$ws.Range("E14").Value = "  -1.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07134"
$ws.Range("E15").Value = "  -1.24%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.38"
$ws.Range("E16").Value = "  +5.35%  "

$ws.Range("E17").Value = "  -2.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009030"
$ws.Range("E18").Value = "  -1.54%  "

$ws.Range("E19").Value = "  -2.36%  "

$ws.Range("E20").Value = "  -0.29%  "

$ws.Range("D21").Value = "27.752.74"
$ws.Range("E21").Value = "  -0.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.270"
$ws.Range("E22").Value = "  -1.06%  "

$ws.Range("E23").Value = "  -1.67%  "

$ws.Range("D24").Value = "2.088.57"
$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.029"
$ws.Range("E25").Value = "  +3.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.75"
$ws.Range("E26").Value = "  -1.78%  "

$ws.Range("E27").Value = "  -1.45%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.132"
$ws.Range("E28").Value = "  +7.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.431"
$ws.Range("E29").Value = "  +1.79%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.67"
$ws.Range("E30").Value = "  +2.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08948"
$ws.Range("E31").Value = "  -1.90%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.239"
$ws.Range("E32").Value = "  +1.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7792"
$ws.Range("E33").Value = "  -0.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.576"
$ws.Range("E34").Value = "  -0.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.915"
$ws.Range("E35").Value = "  -4.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.149"
$ws.Range("E36").Value = "  -1.39%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.011"
$ws.Range("E37").Value = "  -2.49%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05336"
$ws.Range("E38").Value = "  -0.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01972"
$ws.Range("E39").Value = "  -1.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.145"
$ws.Range("E40").Value = "  +3.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.872"
$ws.Range("E41").Value = "  +0.86%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5199"
$ws.Range("E42").Value = "  -0.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1686"
$ws.Range("E43").Value = "  -0.93%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.984"
$ws.Range("E44").Value = "  +2.95%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "110.89"
$ws.Range("E45").Value = "  +0.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.81"
$ws.Range("E46").Value = "  +1.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.722"
$ws.Range("E47").Value = "  -0.40%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4748"
$ws.Range("E48").Value = "  +0.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06504"
$ws.Range("E49").Value = "  +0.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.011"
$ws.Range("E50").Value = "  -2.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.900"
$ws.Range("E51").Value = "  +0.12%  "
